# Update countries & provincias Spain
# Applies the data refresh described by the commit:
#  - bump the "last updated" timestamp
#  - refresh a handful of per-country case counters
#  - reorder "Congo" ahead of "Martinica" in the country list (with Congo's
#    row now carrying the refreshed figures, and Martinica's row carrying
#    the figures that used to belong to the row above it)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados a 21 de Abril de 2020 a las 21:52" -> "...22:22"
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 22:22"

# --- Estados Unidos (row 4): refreshed totals
$ws.Cells.Item(4, 2).Value = 813589
$ws.Cells.Item(4, 3).Value = 20830
$ws.Cells.Item(4, 5).Value = 685956
$ws.Cells.Item(4, 7).Value = 2499
$ws.Cells.Item(4, 8).Value = 45013

# --- Costa Rica (row 95): refreshed totals
$ws.Cells.Item(95, 4).Value = 150
$ws.Cells.Item(95, 5).Value = 513

# --- Niger (row 97): refreshed totals
$ws.Cells.Item(97, 2).Value = 657
$ws.Cells.Item(97, 3).Value = 9
$ws.Cells.Item(97, 4).Value = 127
$ws.Cells.Item(97, 5).Value = 510

# --- Reorder Congo ahead of Martinica (rows 129/130), Congo refreshed
$ws.Cells.Item(129, 1).Value = "Congo"
$ws.Cells.Item(129, 2).Value = 165
$ws.Cells.Item(129, 3).Value = 5
$ws.Cells.Item(129, 4).Value = 16
$ws.Cells.Item(129, 5).Value = 143
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 6

$ws.Cells.Item(130, 1).Value = "Martinica"
$ws.Cells.Item(130, 2).Value = 163
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 73
$ws.Cells.Item(130, 5).Value = 76
$ws.Cells.Item(130, 6).Value = 11
$ws.Cells.Item(130, 7).Value = 2
$ws.Cells.Item(130, 8).Value = 14
